$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp-only rows (column O) ---
$ws.Range("O2").Value = "2022-07-17 20:57:34"
$ws.Range("O3").Value = "2022-07-17 20:57:34"
$ws.Range("O6").Value = "2022-07-17 20:57:34"
$ws.Range("O7").Value = "2022-07-17 20:57:34"
$ws.Range("O8").Value = "2022-07-17 20:57:34"
$ws.Range("O9").Value = "2022-07-17 20:57:34"
$ws.Range("O10").Value = "2022-07-17 20:57:34"
$ws.Range("O11").Value = "2022-07-17 20:57:34"
$ws.Range("O12").Value = "2022-07-17 20:57:34"
$ws.Range("O13").Value = "2022-07-17 20:57:34"
$ws.Range("O14").Value = "2022-07-17 20:57:34"
$ws.Range("O15").Value = "2022-07-17 20:57:34"
$ws.Range("O16").Value = "2022-07-17 20:57:34"
$ws.Range("O17").Value = "2022-07-17 20:57:34"
$ws.Range("O18").Value = "2022-07-17 20:57:34"
$ws.Range("O19").Value = "2022-07-17 20:57:34"
$ws.Range("O23").Value = "2022-07-17 20:57:34"
$ws.Range("O24").Value = "2022-07-17 20:57:34"
$ws.Range("O25").Value = "2022-07-17 20:57:34"
$ws.Range("O26").Value = "2022-07-17 20:57:34"
$ws.Range("O27").Value = "2022-07-17 20:57:34"
$ws.Range("O28").Value = "2022-07-17 20:57:34"
$ws.Range("O31").Value = "2022-07-17 20:57:34"
$ws.Range("O32").Value = "2022-07-17 20:57:34"
$ws.Range("O44").Value = "2022-07-17 20:57:34"
$ws.Range("O45").Value = "2022-07-17 20:57:34"
$ws.Range("O46").Value = "2022-07-17 20:57:34"
$ws.Range("O47").Value = "2022-07-17 20:57:34"
$ws.Range("O48").Value = "2022-07-17 20:57:34"
$ws.Range("O49").Value = "2022-07-17 20:57:34"
$ws.Range("O50").Value = "2022-07-17 20:57:34"
$ws.Range("O51").Value = "2022-07-17 20:57:34"
$ws.Range("O52").Value = "2022-07-17 20:57:34"
$ws.Range("O53").Value = "2022-07-17 20:57:34"
$ws.Range("O54").Value = "2022-07-17 20:57:34"
$ws.Range("O55").Value = "2022-07-17 20:57:34"
$ws.Range("O56").Value = "2022-07-17 20:57:34"
$ws.Range("O57").Value = "2022-07-17 20:57:34"
$ws.Range("O58").Value = "2022-07-17 20:57:34"
$ws.Range("O59").Value = "2022-07-17 20:57:34"
$ws.Range("O60").Value = "2022-07-17 20:57:34"
$ws.Range("O61").Value = "2022-07-17 20:57:34"
$ws.Range("O62").Value = "2022-07-17 20:57:34"
$ws.Range("O63").Value = "2022-07-17 20:57:34"
$ws.Range("O64").Value = "2022-07-17 20:57:34"
$ws.Range("O65").Value = "2022-07-17 20:57:34"
$ws.Range("O66").Value = "2022-07-17 20:57:34"
$ws.Range("O67").Value = "2022-07-17 20:57:34"
$ws.Range("O68").Value = "2022-07-17 20:57:34"
$ws.Range("O69").Value = "2022-07-17 20:57:34"
$ws.Range("O70").Value = "2022-07-17 20:57:34"
$ws.Range("O71").Value = "2022-07-17 20:57:34"
$ws.Range("O72").Value = "2022-07-17 20:57:34"
$ws.Range("O73").Value = "2022-07-17 20:57:34"
$ws.Range("O74").Value = "2022-07-17 20:57:34"
$ws.Range("O75").Value = "2022-07-17 20:57:34"
$ws.Range("O76").Value = "2022-07-17 20:57:34"
$ws.Range("O77").Value = "2022-07-17 20:57:34"
$ws.Range("O78").Value = "2022-07-17 20:57:34"
$ws.Range("O79").Value = "2022-07-17 20:57:34"
$ws.Range("O80").Value = "2022-07-17 20:57:34"
$ws.Range("O81").Value = "2022-07-17 20:57:34"
$ws.Range("O82").Value = "2022-07-17 20:57:34"
$ws.Range("O83").Value = "2022-07-17 20:57:34"
$ws.Range("O84").Value = "2022-07-17 20:57:34"
$ws.Range("O85").Value = "2022-07-17 20:57:34"
$ws.Range("O86").Value = "2022-07-17 20:57:34"
$ws.Range("O87").Value = "2022-07-17 20:57:34"
$ws.Range("O88").Value = "2022-07-17 20:57:34"

# --- Update rows whose content was reordered/changed ---
# Row 4
$ws.Range("A4").Value = "'5920020"
$ws.Range("B4").Value = "BASIC Alkaline Batterien LR6/AA, 12 Stück"
$ws.Range("C4").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/basic-alkaline-batterien-lr6aa-12-stueck/p/5920020"
$ws.Range("D4").Value = "12ST"
$ws.Range("G4").Value = "Coop"
$ws.Range("H4").Value = "'9.95"
$ws.Range("I4").Value = "0.83/1ST"
$ws.Range("J4").Value = "Preis pro 1 Stück"
$ws.Range("K4").Value = "'0.83"
$ws.Range("L4").Value = "1ST"
$ws.Range("M4").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N4").Value = "BASIC Alkaline Batterien LR6/AA, 12 Stück 9.95 Schweizer Franken"
$ws.Range("O4").Value = "2022-07-17 20:57:34"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 0

# Row 5
$ws.Range("A5").Value = "'4488074"
$ws.Range("B5").Value = "Varta Knopfzellen CR2032 2 Stück"
$ws.Range("C5").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2032-2-stueck/p/4488074"
$ws.Range("D5").Value = "2ST"
$ws.Range("G5").Value = "Varta"
$ws.Range("H5").Value = "'8.95"
$ws.Range("I5").Value = "4.48/1ST"
$ws.Range("J5").Value = "Preis pro 1 Stück"
$ws.Range("K5").Value = "'4.48"
$ws.Range("L5").Value = "1ST"
$ws.Range("M5").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N5").Value = "Varta Knopfzellen CR2032 2 Stück 8.95 Schweizer Franken"
$ws.Range("O5").Value = "2022-07-17 20:57:34"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4.5

# Row 20
$ws.Range("A20").Value = "'6761135"
$ws.Range("B20").Value = "Duracell Batterie PLUS 9V/6LR61 1 Stück"
$ws.Range("C20").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-plus-9v6lr61-1-stueck/p/6761135"
$ws.Range("D20").Value = "1ST"
$ws.Range("G20").Value = "Duracell"
$ws.Range("H20").Value = "'9.95"
$ws.Range("I20").Value = "9.95/1ST"
$ws.Range("J20").Value = "Preis pro 1 Stück"
$ws.Range("K20").Value = "'9.95"
$ws.Range("L20").Value = "1ST"
$ws.Range("M20").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N20").Value = "Duracell Batterie PLUS 9V/6LR61 1 Stück 9.95 Schweizer Franken"
$ws.Range("O20").Value = "2022-07-17 20:57:34"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = 0

# Row 21
$ws.Range("A21").Value = "'3494233"
$ws.Range("B21").Value = "Varta Electronics CR2032 1er Bli"
$ws.Range("C21").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2032-1er-bli/p/3494233"
$ws.Range("D21").Value = "1ST"
$ws.Range("G21").Value = "Varta"
$ws.Range("H21").Value = "'4.95"
$ws.Range("I21").Value = "4.95/1ST"
$ws.Range("J21").Value = "Preis pro 1 Stück"
$ws.Range("K21").Value = "'4.95"
$ws.Range("L21").Value = "1ST"
$ws.Range("M21").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N21").Value = "Varta Electronics CR2032 1er Bli 4.95 Schweizer Franken"
$ws.Range("O21").Value = "2022-07-17 20:57:34"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 4.5

# Row 22
$ws.Range("A22").Value = "'6753557"
$ws.Range("B22").Value = "Duracell Batterien PLUS AAA/LR03 4 Stück"
$ws.Range("C22").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-4-stueck/p/6753557"
$ws.Range("D22").Value = "4ST"
$ws.Range("G22").Value = "Duracell"
$ws.Range("H22").Value = "'9.95"
$ws.Range("I22").Value = "2.49/1ST"
$ws.Range("J22").Value = "Preis pro 1 Stück"
$ws.Range("K22").Value = "'2.49"
$ws.Range("L22").Value = "1ST"
$ws.Range("M22").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N22").Value = "Duracell Batterien PLUS AAA/LR03 4 Stück 9.95 Schweizer Franken"
$ws.Range("O22").Value = "2022-07-17 20:57:34"
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = 0

# Row 29
$ws.Range("A29").Value = "'6508223"
$ws.Range("B29").Value = "satrap Venti WS Standventilator"
$ws.Range("C29").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-venti-ws-standventilator/p/6508223"
$ws.Range("D29").Value = ""
$ws.Range("G29").Value = "satrap"
$ws.Range("H29").Value = "'29.95"
$ws.Range("I29").Value = ""
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N29").Value = "satrap Venti WS Standventilator 50% Aktion 29.95 Schweizer Franken statt 59.95 Schweizer Franken"
$ws.Range("O29").Value = "2022-07-17 20:57:34"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = 0

# Row 30
$ws.Range("A30").Value = "'3494130"
$ws.Range("B30").Value = "Varta Longlife Power AA 4er Bli"
$ws.Range("C30").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-4er-bli/p/3494130"
$ws.Range("D30").Value = "4ST"
$ws.Range("G30").Value = "Varta"
$ws.Range("H30").Value = "'8.95"
$ws.Range("I30").Value = "2.24/1ST"
$ws.Range("J30").Value = "Preis pro 1 Stück"
$ws.Range("K30").Value = "'2.24"
$ws.Range("L30").Value = "1ST"
$ws.Range("M30").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N30").Value = "Varta Longlife Power AA 4er Bli 8.95 Schweizer Franken"
$ws.Range("O30").Value = "2022-07-17 20:57:34"
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 4

# Row 33
$ws.Range("A33").Value = "'6999781"
$ws.Range("B33").Value = "Varta Batterien Longlife Power AAA/LR03 2x12 Stück"
$ws.Range("C33").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-batterien-longlife-power-aaalr03-2x12-stueck/p/6999781"
$ws.Range("D33").Value = "24ST"
$ws.Range("G33").Value = "Varta"
$ws.Range("H33").Value = "'20.85"
$ws.Range("I33").Value = "0.87/1ST"
$ws.Range("J33").Value = "Preis pro 1 Stück"
$ws.Range("K33").Value = "'0.87"
$ws.Range("L33").Value = "1ST"
$ws.Range("M33").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N33").Value = "Varta Batterien Longlife Power AAA/LR03 2x12 Stück 50% Aktion 20.85 Schweizer Franken statt 41.70 Schweizer Franken"
$ws.Range("O33").Value = "2022-07-17 20:57:34"
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = 0

# Row 34
$ws.Range("A34").Value = "'4905484"
$ws.Range("B34").Value = "Alkaline Batterie LR20/D 2 Stück"
$ws.Range("C34").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr20d-2-stueck/p/4905484"
$ws.Range("D34").Value = "2ST"
$ws.Range("G34").Value = "Coop"
$ws.Range("H34").Value = "'5.95"
$ws.Range("I34").Value = "2.98/1ST"
$ws.Range("J34").Value = "Preis pro 1 Stück"
$ws.Range("K34").Value = "'2.98"
$ws.Range("L34").Value = "1ST"
$ws.Range("M34").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N34").Value = "Alkaline Batterie LR20/D 2 Stück 5.95 Schweizer Franken"
$ws.Range("O34").Value = "2022-07-17 20:57:34"
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 5

# Row 35
$ws.Range("A35").Value = "'6999749"
$ws.Range("B35").Value = "Varta Batterien Longlife Power AA/LR6 2x12 Stück"
$ws.Range("C35").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-batterien-longlife-power-aalr6-2x12-stueck/p/6999749"
$ws.Range("D35").Value = "24ST"
$ws.Range("G35").Value = "Varta"
$ws.Range("H35").Value = "'20.85"
$ws.Range("I35").Value = "0.87/1ST"
$ws.Range("J35").Value = "Preis pro 1 Stück"
$ws.Range("K35").Value = "'0.87"
$ws.Range("L35").Value = "1ST"
$ws.Range("M35").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N35").Value = "Varta Batterien Longlife Power AA/LR6 2x12 Stück 50% Aktion 20.85 Schweizer Franken statt 41.70 Schweizer Franken"
$ws.Range("O35").Value = "2022-07-17 20:57:34"
$ws.Range("E35").Value = ""
$ws.Range("F35").Value = 0

# Row 36
$ws.Range("A36").Value = "'6753555"
$ws.Range("B36").Value = "Duracell Batterien Optimum AAA/LR03 4 Stück"
$ws.Range("C36").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-optimum-aaalr03-4-stueck/p/6753555"
$ws.Range("D36").Value = "4ST"
$ws.Range("G36").Value = "Duracell"
$ws.Range("H36").Value = "'11.95"
$ws.Range("I36").Value = "2.99/1ST"
$ws.Range("J36").Value = "Preis pro 1 Stück"
$ws.Range("K36").Value = "'2.99"
$ws.Range("L36").Value = "1ST"
$ws.Range("M36").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N36").Value = "Duracell Batterien Optimum AAA/LR03 4 Stück 11.95 Schweizer Franken"
$ws.Range("O36").Value = "2022-07-17 20:57:34"
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = 0

# Row 37
$ws.Range("A37").Value = "'3494138"
$ws.Range("B37").Value = "Varta Longlife Power Batterien AAA/LR03 8 Stück"
$ws.Range("C37").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-batterien-aaalr03-8-stueck/p/3494138"
$ws.Range("D37").Value = "8ST"
$ws.Range("G37").Value = "Varta"
$ws.Range("H37").Value = "'14.95"
$ws.Range("I37").Value = "1.87/1ST"
$ws.Range("J37").Value = "Preis pro 1 Stück"
$ws.Range("K37").Value = "'1.87"
$ws.Range("L37").Value = "1ST"
$ws.Range("M37").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N37").Value = "Varta Longlife Power Batterien AAA/LR03 8 Stück 14.95 Schweizer Franken"
$ws.Range("O37").Value = "2022-07-17 20:57:34"
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 4

# Row 38
$ws.Range("A38").Value = "'4942597"
$ws.Range("B38").Value = "Skross Adapter Europa-Schweiz"
$ws.Range("C38").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-europa-schweiz/p/4942597"
$ws.Range("D38").Value = ""
$ws.Range("G38").Value = "Skross"
$ws.Range("H38").Value = "'14.95"
$ws.Range("I38").Value = ""
$ws.Range("J38").Value = ""
$ws.Range("K38").Value = ""
$ws.Range("L38").Value = ""
$ws.Range("M38").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N38").Value = "Skross Adapter Europa-Schweiz 14.95 Schweizer Franken"
$ws.Range("O38").Value = "2022-07-17 20:57:34"
$ws.Range("E38").Value = ""
$ws.Range("F38").Value = 0

# Row 39
$ws.Range("A39").Value = "'3494909"
$ws.Range("B39").Value = "Varta Electronics V23GA 1er Bli"
$ws.Range("C39").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v23ga-1er-bli/p/3494909"
$ws.Range("D39").Value = "1ST"
$ws.Range("G39").Value = "Varta"
$ws.Range("H39").Value = "'3.95"
$ws.Range("I39").Value = "3.95/1ST"
$ws.Range("J39").Value = "Preis pro 1 Stück"
$ws.Range("K39").Value = "'3.95"
$ws.Range("L39").Value = "1ST"
$ws.Range("M39").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N39").Value = "Varta Electronics V23GA 1er Bli 3.95 Schweizer Franken"
$ws.Range("O39").Value = "2022-07-17 20:57:34"
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 5

# Row 40
$ws.Range("A40").Value = "'4096751"
$ws.Range("B40").Value = "Varta Longlife Power Batterien AA/LR6 6 Stück"
$ws.Range("C40").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-batterien-aalr6-6-stueck/p/4096751"
$ws.Range("D40").Value = "6ST"
$ws.Range("G40").Value = "Varta"
$ws.Range("H40").Value = "'12.95"
$ws.Range("I40").Value = "2.16/1ST"
$ws.Range("J40").Value = "Preis pro 1 Stück"
$ws.Range("K40").Value = "'2.16"
$ws.Range("L40").Value = "1ST"
$ws.Range("M40").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N40").Value = "Varta Longlife Power Batterien AA/LR6 6 Stück 12.95 Schweizer Franken"
$ws.Range("O40").Value = "2022-07-17 20:57:34"
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = 0

# Row 41
$ws.Range("A41").Value = "'6848736"
$ws.Range("B41").Value = "SONY WF-C500B (In-Ear, Bluetooth 5.0, Schwarz)"
$ws.Range("C41").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/sony-wf-c500b-in-ear-bluetooth-50-schwarz/p/6848736"
$ws.Range("D41").Value = ""
$ws.Range("G41").Value = "Sony"
$ws.Range("H41").Value = "'39.95"
$ws.Range("I41").Value = ""
$ws.Range("J41").Value = ""
$ws.Range("K41").Value = ""
$ws.Range("L41").Value = ""
$ws.Range("M41").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete']"
$ws.Range("N41").Value = "SONY WF-C500B (In-Ear, Bluetooth 5.0, Schwarz) 50% Aktion 39.95 Schweizer Franken statt 79.95 Schweizer Franken"
$ws.Range("O41").Value = "2022-07-17 20:57:34"
$ws.Range("E41").Value = ""
$ws.Range("F41").Value = 0

# Row 42
$ws.Range("A42").Value = "'5751576"
$ws.Range("B42").Value = "satrap Toasty 1 Toaster"
$ws.Range("C42").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-toasty-1-toaster/p/5751576"
$ws.Range("D42").Value = ""
$ws.Range("G42").Value = "satrap"
$ws.Range("H42").Value = "'29.95"
$ws.Range("I42").Value = ""
$ws.Range("J42").Value = ""
$ws.Range("K42").Value = ""
$ws.Range("L42").Value = ""
$ws.Range("M42").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N42").Value = "satrap Toasty 1 Toaster 29.95 Schweizer Franken"
$ws.Range("O42").Value = "2022-07-17 20:57:34"
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 5

# Row 43
$ws.Range("A43").Value = "'6119284"
$ws.Range("B43").Value = "satrap Aqua SA10 Wasserkocher"
$ws.Range("C43").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-aqua-sa10-wasserkocher/p/6119284"
$ws.Range("D43").Value = ""
$ws.Range("G43").Value = "satrap"
$ws.Range("H43").Value = "'49.95"
$ws.Range("I43").Value = ""
$ws.Range("J43").Value = ""
$ws.Range("K43").Value = ""
$ws.Range("L43").Value = ""
$ws.Range("M43").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N43").Value = "satrap Aqua SA10 Wasserkocher 49.95 Schweizer Franken"
$ws.Range("O43").Value = "2022-07-17 20:57:34"
$ws.Range("E43").Value = 3
$ws.Range("F43").Value = 2.5

Write-Host "Edit complete"
